$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ReifenDaten")

# Update the tire-data area formulas (rows 2-258):
#   Y column: (0.32*(K-V))/200  ->  (32*(L-V))/200
#   Z column: 2.54*Y*SQRT((K-Y)*(M-Y))  ->  2.54*Y*SQRT((L-Y)*(N-Y))
$ws.Range("Y2:Y258").FormulaR1C1 = "=(32*(RC[-13]-RC[-3]))/200"
$ws.Range("Z2:Z258").FormulaR1C1 = "=2.54*RC[-1]*SQRT((RC[-14]-RC[-1])*(RC[-12]-RC[-1]))"

# Update the view so the new selection matches the saved state
$ws.Application.Goto($ws.Range("L193"), $true)
$ws.Range("Y202").Select()
